$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1): update F3, F5, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3215
$ws1.Range("F5").Value = 24
$ws1.Range("F6").Value = 138

# Sheet "全部类型" (index 4): update F7, F9, F11
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3215
$ws4.Range("F9").Value = 24
$ws4.Range("F11").Value = 138
